$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@("ECs","Gdf11","Acvr2b","ECs",3,1,1.296436,3.889308000000001,0.1944674516316147,0.1944674516316147,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,1.996077998230667,17.964701984076,0.06946549071113942,0.06946549071113942)
  ,@("ECs","Gdf11","Acvr2b","FAPs",3,1,1.296436,3.889308000000001,0.1944674516316147,0.1944674516316147,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,1.883555996337334,16.952003967036,0.06554961363406694,0.06554961363406694)
  ,@("ECs","Gdf11","Acvr2b","sCs",3,1,1.296436,3.889308000000001,0.1944674516316147,0.1944674516316147,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,1.708352178134667,15.375169603212,0.0594523472864084,0.0594523472864084)
  ,@("FAPs","Gdf11","Acvr2b","ECs",3,1,1.837275666666667,5.511827,0.2755942575194169,0.2755942575194169,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,2.828790264168778,25.459112377519,0.09844470205751446,0.09844470205751446)
  ,@("FAPs","Gdf11","Acvr2b","FAPs",3,1,1.837275666666667,5.511827,0.2755942575194169,0.2755942575194169,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,2.669326984806556,24.023942863259,0.09289522204665156,0.09289522204665156)
  ,@("FAPs","Gdf11","Acvr2b","sCs",3,1,1.837275666666667,5.511827,0.2755942575194169,0.2755942575194169,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,2.421032651811445,21.789293866303,0.08425433341525086,0.08425433341525086)
  ,@("sCs","Gdf11","Acvr2b","ECs",3,1,3.532884666666666,10.598654,0.5299382908489685,0.5299382908489685,2,0.6666666666666666,1.539665666666667,4.618997,0.3572088291809875,0.3572088291809875,5.439461225559778,48.955151030038,0.1892986364123337,0.1892986364123337)
  ,@("sCs","Gdf11","Acvr2b","FAPs",3,1,3.532884666666666,10.598654,0.5299382908489685,0.5299382908489685,3,1,1.452872333333333,4.358617,0.3370724153789985,0.3370724153789985,5.132830389057555,46.19547350151799,0.17862757969828,0.17862757969828)
  ,@("sCs","Gdf11","Acvr2b","sCs",3,1,3.532884666666666,10.598654,0.5299382908489685,0.5299382908489685,3,1,1.317729666666667,3.953189,0.3057187554400141,0.3057187554400141,4.655386934178444,41.898482407606,0.1620120747383549,0.1620120747383549)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = 2 + $i
  $row = $rows[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $c = 1 + $j
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}
